$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "Actual" effort column (J3:J15) with the recorded actuals
$ws.Range("J3").Value  = 3
$ws.Range("J4").Value  = 2
$ws.Range("J5").Value  = 40
$ws.Range("J6").Value  = 2
$ws.Range("J7").Value  = 1
$ws.Range("J8").Value  = 1
$ws.Range("J9").Value  = 3
$ws.Range("J10").Value = 2
$ws.Range("J11").Value = 2
$ws.Range("J12").Value = 5
$ws.Range("J13").Value = 21
$ws.Range("J14").Value = 6
$ws.Range("J15").Value = 21

# Add totals + variance row underneath the table
$ws.Range("I16").Formula = "=SUM(I3:I15)"
$ws.Range("J16").Formula = "=SUM(J3:J15)"
$ws.Range("K16").Formula = "=1-J16/I16"
$ws.Range("K16").Style = "Percent"

# Move the viewport / active selection like the author left it
$ws.Range("L8").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 20

Write-Output "done"
